$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = [double]"35.160799"
    "H2" = [double]"105.482397"
    "I2" = [double]"0.2238945559395223"
    "J2" = [double]"0.2238945559395223"
    "M2" = [double]"0.1825283333333333"
    "N2" = [double]"0.547585"
    "O2" = [double]"0.001028331058213739"
    "P2" = [double]"0.001028331058213739"
    "Q2" = [double]"6.417842040138332"
    "R2" = [double]"57.760578361245"
    "S2" = [double]"0.000230237725637584"
    "T2" = [double]"0.0002302377256375841"
    "G3" = [double]"35.160799"
    "H3" = [double]"105.482397"
    "I3" = [double]"0.2238945559395223"
    "J3" = [double]"0.2238945559395223"
    "O3" = [double]"0.0001759459539160193"
    "P3" = [double]"0.0001759459539160193"
    "Q3" = [double]"1.098083473036333"
    "R3" = [double]"9.882751257327"
    "S3" = [double]"3.939334122138278E-05"
    "T3" = [double]"3.939334122138278E-05"
    "G4" = [double]"35.160799"
    "H4" = [double]"105.482397"
    "I4" = [double]"0.2238945559395223"
    "J4" = [double]"0.2238945559395223"
    "M4" = [double]"103.239782"
    "N4" = [double]"309.719346"
    "O4" = [double]"0.5816339432625932"
    "P4" = [double]"0.5816339432625932"
    "Q4" = [double]"3629.993223705817"
    "R4" = [double]"32669.93901335236"
    "S4" = [double]"0.1302246734461316"
    "T4" = [double]"0.1302246734461316"
    "G5" = [double]"35.160799"
    "H5" = [double]"105.482397"
    "I5" = [double]"0.2238945559395223"
    "J5" = [double]"0.2238945559395223"
    "M5" = [double]"0.04852733333333333"
    "N5" = [double]"0.145582"
    "O5" = [double]"0.0002733940705404138"
    "P5" = [double]"0.0002733940705404139"
    "Q5" = [double]"1.706259813339333"
    "R5" = [double]"15.356338320054"
    "S5" = [double]"6.121144402014438E-05"
    "T5" = [double]"6.121144402014439E-05"
    "G6" = [double]"35.160799"
    "H6" = [double]"105.482397"
    "I6" = [double]"0.2238945559395223"
    "J6" = [double]"0.2238945559395223"
    "M6" = [double]"73.99751433333334"
    "N6" = [double]"221.992543"
    "O6" = [double]"0.4168883856547366"
    "P6" = [double]"0.4168883856547366"
    "Q6" = [double]"2601.811727973953"
    "R6" = [double]"23416.30555176557"
    "S6" = [double]"0.09333903998251156"
    "T6" = [double]"0.09333903998251156"
    "I7" = [double]"0.2050965007332699"
    "J7" = [double]"0.2050965007332699"
    "M7" = [double]"0.1825283333333333"
    "N7" = [double]"0.547585"
    "O7" = [double]"0.001028331058213739"
    "P7" = [double]"0.001028331058213739"
    "Q7" = [double]"5.879003797871666"
    "R7" = [double]"52.911034180845"
    "S7" = [double]"0.0002109071016349782"
    "T7" = [double]"0.0002109071016349783"
    "I8" = [double]"0.2050965007332699"
    "J8" = [double]"0.2050965007332699"
    "O8" = [double]"0.0001759459539160193"
    "P8" = [double]"0.0001759459539160193"
    "S8" = [double]"3.608589946635271E-05"
    "T8" = [double]"3.608589946635271E-05"
    "I9" = [double]"0.2050965007332699"
    "J9" = [double]"0.2050965007332699"
    "M9" = [double]"103.239782"
    "N9" = [double]"309.719346"
    "O9" = [double]"0.5816339432625932"
    "P9" = [double]"0.5816339432625932"
    "Q9" = [double]"3325.221128059258"
    "R9" = [double]"29926.99015253332"
    "S9" = [double]"0.1192910864708511"
    "T9" = [double]"0.1192910864708511"
    "I10" = [double]"0.2050965007332699"
    "J10" = [double]"0.2050965007332699"
    "M10" = [double]"0.04852733333333333"
    "N10" = [double]"0.145582"
    "O10" = [double]"0.0002733940705404138"
    "P10" = [double]"0.0002733940705404139"
    "Q10" = [double]"1.563003243152667"
    "R10" = [double]"14.067029188374"
    "S10" = [double]"5.607216718906362E-05"
    "T10" = [double]"5.607216718906363E-05"
    "I11" = [double]"0.2050965007332699"
    "J11" = [double]"0.2050965007332699"
    "M11" = [double]"73.99751433333334"
    "N11" = [double]"221.992543"
    "O11" = [double]"0.4168883856547366"
    "P11" = [double]"0.4168883856547366"
    "Q11" = [double]"2383.365145860806"
    "R11" = [double]"21450.28631274725"
    "S11" = [double]"0.08550234909412839"
    "T11" = [double]"0.08550234909412839"
    "G12" = [double]"51.53356533333334"
    "H12" = [double]"154.600696"
    "I12" = [double]"0.3281519491717758"
    "J12" = [double]"0.3281519491717758"
    "M12" = [double]"0.1825283333333333"
    "N12" = [double]"0.547585"
    "O12" = [double]"0.001028331058213739"
    "P12" = [double]"0.001028331058213739"
    "Q12" = [double]"9.406335791017778"
    "R12" = [double]"84.65702211916"
    "S12" = [double]"0.0003374488411467132"
    "T12" = [double]"0.0003374488411467133"
    "G13" = [double]"51.53356533333334"
    "H13" = [double]"154.600696"
    "I13" = [double]"0.3281519491717758"
    "J13" = [double]"0.3281519491717758"
    "O13" = [double]"0.0001759459539160193"
    "P13" = [double]"0.0001759459539160193"
    "Q13" = [double]"1.609410423215111"
    "R13" = [double]"14.484693808936"
    "S13" = [double]"5.773700772642916E-05"
    "T13" = [double]"5.773700772642916E-05"
    "G14" = [double]"51.53356533333334"
    "H14" = [double]"154.600696"
    "I14" = [double]"0.3281519491717758"
    "J14" = [double]"0.3281519491717758"
    "M14" = [double]"103.239782"
    "N14" = [double]"309.719346"
    "O14" = [double]"0.5816339432625932"
    "P14" = [double]"0.5816339432625932"
    "Q14" = [double]"5320.314050696091"
    "R14" = [double]"47882.82645626481"
    "S14" = [double]"0.190864312186086"
    "T14" = [double]"0.190864312186086"
    "G15" = [double]"51.53356533333334"
    "H15" = [double]"154.600696"
    "I15" = [double]"0.3281519491717758"
    "J15" = [double]"0.3281519491717758"
    "M15" = [double]"0.04852733333333333"
    "N15" = [double]"0.145582"
    "O15" = [double]"0.0002733940705404138"
    "P15" = [double]"0.0002733940705404139"
    "Q15" = [double]"2.500786502785778"
    "R15" = [double]"22.507078525072"
    "S15" = [double]"8.971479713984277E-05"
    "T15" = [double]"8.971479713984278E-05"
    "G16" = [double]"51.53356533333334"
    "H16" = [double]"154.600696"
    "I16" = [double]"0.3281519491717758"
    "J16" = [double]"0.3281519491717758"
    "M16" = [double]"73.99751433333334"
    "N16" = [double]"221.992543"
    "O16" = [double]"0.4168883856547366"
    "P16" = [double]"0.4168883856547366"
    "Q16" = [double]"3813.355739401104"
    "R16" = [double]"34320.20165460993"
    "S16" = [double]"0.1368027363396768"
    "T16" = [double]"0.1368027363396768"
    "G17" = [double]"0.5955593333333333"
    "H17" = [double]"1.786678"
    "I17" = [double]"0.003792362411113143"
    "J17" = [double]"0.003792362411113143"
    "M17" = [double]"0.1825283333333333"
    "N17" = [double]"0.547585"
    "O17" = [double]"0.001028331058213739"
    "P17" = [double]"0.001028331058213739"
    "Q17" = [double]"0.1087064525144444"
    "R17" = [double]"0.97835807263"
    "S17" = [double]"3.899804051349984E-06"
    "T17" = [double]"3.899804051349985E-06"
    "G18" = [double]"0.5955593333333333"
    "H18" = [double]"1.786678"
    "I18" = [double]"0.003792362411113143"
    "J18" = [double]"0.003792362411113143"
    "O18" = [double]"0.0001759459539160193"
    "P18" = [double]"0.0001759459539160193"
    "Q18" = [double]"0.01859951649977778"
    "R18" = [double]"0.167395648498"
    "S18" = [double]"6.672508220185568E-07"
    "T18" = [double]"6.672508220185568E-07"
    "G19" = [double]"0.5955593333333333"
    "H19" = [double]"1.786678"
    "I19" = [double]"0.003792362411113143"
    "J19" = [double]"0.003792362411113143"
    "M19" = [double]"103.239782"
    "N19" = [double]"309.719346"
    "O19" = [double]"0.5816339432625932"
    "P19" = [double]"0.5816339432625932"
    "Q19" = [double]"61.48541574139866"
    "R19" = [double]"553.3687416725879"
    "S19" = [double]"0.002205766703456573"
    "T19" = [double]"0.002205766703456573"
    "G20" = [double]"0.5955593333333333"
    "H20" = [double]"1.786678"
    "I20" = [double]"0.003792362411113143"
    "J20" = [double]"0.003792362411113143"
    "M20" = [double]"0.04852733333333333"
    "N20" = [double]"0.145582"
    "O20" = [double]"0.0002733940705404138"
    "P20" = [double]"0.0002733940705404139"
    "Q20" = [double]"0.02890090628844444"
    "R20" = [double]"0.260108156596"
    "S20" = [double]"1.036809396538681E-06"
    "T20" = [double]"1.036809396538681E-06"
    "G21" = [double]"0.5955593333333333"
    "H21" = [double]"1.786678"
    "I21" = [double]"0.003792362411113143"
    "J21" = [double]"0.003792362411113143"
    "M21" = [double]"73.99751433333334"
    "N21" = [double]"221.992543"
    "O21" = [double]"0.4168883856547366"
    "P21" = [double]"0.4168883856547366"
    "Q21" = [double]"44.06991030468378"
    "R21" = [double]"396.629192742154"
    "S21" = [double]"0.001580991843386663"
    "T21" = [double]"0.001580991843386663"
    "G22" = [double]"37.54313466666667"
    "H22" = [double]"112.629404"
    "I22" = [double]"0.2390646317443189"
    "J22" = [double]"0.2390646317443189"
    "M22" = [double]"0.1825283333333333"
    "N22" = [double]"0.547585"
    "O22" = [double]"0.001028331058213739"
    "P22" = [double]"0.001028331058213739"
    "Q22" = [double]"6.852685798815555"
    "R22" = [double]"61.67417218934"
    "S22" = [double]"0.0002458375857431132"
    "T22" = [double]"0.0002458375857431133"
    "G23" = [double]"37.54313466666667"
    "H23" = [double]"112.629404"
    "I23" = [double]"0.2390646317443189"
    "J23" = [double]"0.2390646317443189"
    "O23" = [double]"0.0001759459539160193"
    "P23" = [double]"0.0001759459539160193"
    "Q23" = [double]"1.172484610018222"
    "R23" = [double]"10.552361490164"
    "S23" = [double]"4.206245467983605E-05"
    "T23" = [double]"4.206245467983605E-05"
    "G24" = [double]"37.54313466666667"
    "H24" = [double]"112.629404"
    "I24" = [double]"0.2390646317443189"
    "J24" = [double]"0.2390646317443189"
    "M24" = [double]"103.239782"
    "N24" = [double]"309.719346"
    "O24" = [double]"0.5816339432625932"
    "P24" = [double]"0.5816339432625932"
    "Q24" = [double]"3875.945038583309"
    "R24" = [double]"34883.50534724978"
    "S24" = [double]"0.1390481044560679"
    "T24" = [double]"0.1390481044560679"
    "G25" = [double]"37.54313466666667"
    "H25" = [double]"112.629404"
    "I25" = [double]"0.2390646317443189"
    "J25" = [double]"0.2390646317443189"
    "M25" = [double]"0.04852733333333333"
    "N25" = [double]"0.145582"
    "O25" = [double]"0.0002733940705404138"
    "P25" = [double]"0.0002733940705404139"
    "Q25" = [double]"1.821868210347555"
    "R25" = [double]"16.396813893128"
    "S25" = [double]"6.535885279482439E-05"
    "T25" = [double]"6.53588527948244E-05"
    "G26" = [double]"37.54313466666667"
    "H26" = [double]"112.629404"
    "I26" = [double]"0.2390646317443189"
    "J26" = [double]"0.2390646317443189"
    "M26" = [double]"73.99751433333334"
    "N26" = [double]"221.992543"
    "O26" = [double]"0.4168883856547366"
    "P26" = [double]"0.4168883856547366"
    "Q26" = [double]"2778.098645614931"
    "R26" = [double]"25002.88781053437"
    "S26" = [double]"0.09966326839503323"
    "T26" = [double]"0.09966326839503323"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
